# Apply the changes from the commit "Add UML sequence diagrams content (#3)"
# to the use-case "inclusion" diagram deck.
#
# The user-visible content edit captured by the diff is the removal of the
# standalone slide title textbox ("Use case diagram sample") that used to
# sit above the diagram -- the slide keeps all of its diagram shapes but no
# longer carries its own title placeholder/text.
#
# The deck's custom starting slide number was also reset back to the
# default as part of the same save.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate and remove the title placeholder shape. Cut() (rather than
# Delete()) is used because Delete()-ing a title placeholder just leaves an
# empty "click to add title" placeholder behind (the layout still supplies
# one), whereas Cut() removes it from the slide's shape tree completely --
# matching the target deck, which has no title shape left on this slide.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    $isTitle = $false
    if ($sh.Type -eq 14) {
        if ($sh.PlaceholderFormat.Type -eq 1 -or $sh.PlaceholderFormat.Type -eq 13) {
            $isTitle = $true
        }
    }
    if ($sh.Name -eq "Title 2") {
        $isTitle = $true
    }
    if ($isTitle) {
        $sh.Cut()
    }
}

# The deck no longer uses a custom starting slide number -- restore the
# default (1) instead of the previous offset (37).
$p.PageSetup.FirstSlideNumber = 1
